# graviton_master.xlsx — "further detailing graviton_master spreadsheet"
#
# 1. Insert a new empty "phase_noise" sheet right after "rx_chain".
# 2. Insert a new "power_budget" sheet (with a small supply table) right
#    after "tx_chain".
# 3. Update a few input cells on "rx_chain" (E13, E20) that ripple through
#    the existing gain/NF cascade formulas, and give G8 an explicit
#    =ABS(E8) formula (matching its siblings G9:G22) instead of a bare
#    literal.
# 4. Re-point the active tab/selection: rx_chain becomes the active sheet
#    (selection on E28) instead of tx_chain.

$wb = $excel.ActiveWorkbook

$rxSeed = $wb.Worksheets.Item("rx_chain")
$txSeed = $wb.Worksheets.Item("tx_chain")

# --- new sheets -----------------------------------------------------
# Create power_budget first and phase_noise second so the engine hands
# out sheetId=3 to power_budget and sheetId=4 to phase_noise, matching
# the target sheet order: rx_chain, phase_noise, tx_chain, power_budget.
$powerBudget = $wb.Worksheets.Add($null, $txSeed)
$powerBudget.Name = "power_budget"

$phaseNoise = $wb.Worksheets.Add($null, $rxSeed)
$phaseNoise.Name = "phase_noise"

# Worksheet handles captured before an Add() can end up re-pointed at a
# different sheet once tab positions shift, so re-resolve every handle
# by name now that the final sheet count/order is settled, and only do
# the real work (cell writes / selections) from here on.
$rx = $wb.Worksheets.Item("rx_chain")
$tx = $wb.Worksheets.Item("tx_chain")
$powerBudget = $wb.Worksheets.Item("power_budget")
$phaseNoise = $wb.Worksheets.Item("phase_noise")

# --- power_budget contents -------------------------------------------
$powerBudget.Range("G9").Value = "Supply"
$powerBudget.Range("H9").Value = "Source"
$powerBudget.Range("I9").Value = "Imax"
$powerBudget.Range("J9").Value = "Vmin"
$powerBudget.Range("K9").Value = "Vmax"

$powerBudget.Range("G10").Value = "PP5V"
$powerBudget.Range("H10").Value = "VIN48V"
$powerBudget.Range("I10").Value = "6A"
$powerBudget.Range("J10").Formula = "=4.98*0.98"
$powerBudget.Range("K10").Formula = "=4.98*1.02"

# Resting selection for the new sheet (matches the authored workbook).
$powerBudget.Range("F39").Select()

# phase_noise is left blank (placeholder sheet for future content);
# default A1 selection is fine, nothing else to do there.

# --- rx_chain cascade tweaks ------------------------------------------
$rx.Range("E13").Value = -8
$rx.Range("E20").Value = -4
$rx.Range("G8").Formula = "=ABS(E8)"

# --- view / active-tab bookkeeping ------------------------------------
# tx_chain keeps its own resting selection (L23), just no longer active.
$tx.Range("L23").Select()

# rx_chain becomes the active sheet/tab again, selection parked on E28.
# This must run last so it "wins" as the final active tab.
$rx.Activate()
$rx.Range("E28").Select()
